$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Locate the "ACME Movie Database" paragraph (currently styled BodyText2) ---
$pAcme = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "^ACME Movie Database") {
        $pAcme = $cand
        break
    }
}

# 1) Re-style the title paragraph: Heading1, centered, numbering explicitly suppressed
$pAcme.Style = "Heading1"
$pAcme.Range.ListFormat.RemoveNumbers()
$pAcme.Alignment = 1

# 2) Insert the blank "spacer" heading paragraph (indented, centered) right after it
$pAcme.Range.InsertParagraphAfter()
$pSpacer = $pAcme.Next()
$pSpacer.Range.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr><w:jc w:val=`"center`"/></w:pPr></w:p>")
# (set the indent as a follow-up property write - doing it inline above gets
#  silently optimized away because it matches the Heading1 list level's indent)
$pSpacer.LeftIndent = 21.6

# 3) Insert the "Test Documentation" heading paragraph after the spacer
$pSpacer.Range.InsertParagraphAfter()
$pTestDoc = $pSpacer.Next()
$pTestDoc.Range.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr><w:jc w:val=`"center`"/></w:pPr><w:r><w:t xml:space=`"preserve`">Test </w:t></w:r><w:r><w:t>Documentation</w:t></w:r></w:p>")

# 4) Insert the "Sprint 2" heading paragraph (larger type size) after that
$pTestDoc.Range.InsertParagraphAfter()
$pSprint = $pTestDoc.Next()
$pSprint.Range.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val=`"Heading1`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr><w:jc w:val=`"center`"/><w:rPr><w:sz w:val=`"36`"/><w:szCs w:val=`"16`"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=`"36`"/><w:szCs w:val=`"16`"/></w:rPr><w:t>Sprint 2</w:t></w:r></w:p>")

# --- Remove the old (now-duplicated) "Test Documentation" / blank / "Sprint 2" paragraphs ---
# They immediately follow the pre-existing blank paragraph that sits after $pSprint.
$pOldBlank = $pSprint.Next()
$pOldTestDoc = $pOldBlank.Next()
$pOldTestDoc.Range.Delete()              # old "Test Documentation" paragraph
$pOldBlank.Next().Range.Delete()          # old blank BodyText2 paragraph
$pOldBlank.Next().Range.Delete()          # old "Sprint 2" paragraph

# --- Add 15 plain empty paragraphs where the removed ones used to be ---
$anchor = $pOldBlank
for ($n = 0; $n -lt 15; $n++) {
    $anchor.Range.InsertParagraphAfter()
    $newP = $anchor.Next()
    $newP.Range.InsertXML("<w:p $wns></w:p>")
    $anchor = $newP
}
